# Column G ("Recorded By") contains comma-separated lists of two
# names/emails (e.g. "dnasr281@gmail.com, System"). This swaps the order
# of the two entries whenever the first entry is not "System" (cells
# that already start with "System" are left untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1
if ($lastRow -lt 2) { $lastRow = 157 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -eq $val -or $val -eq "") { continue }

    $parts = $val -split ", "
    if ($parts.Count -eq 2 -and $parts[0] -ne "System") {
        $cell.Value2 = "$($parts[1]), $($parts[0])"
    }
}
